$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19.94943145514285
$ws.Range("C2").Value = 9.864324008596844
$ws.Range("D2").Value = 7.119280444993096
$ws.Range("E2").Value = 9.646071058415391
$ws.Range("F2").Value = 37.38304194697183
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 29.55346602580135
$ws.Range("L2").Value = 10.27488303849287
$ws.Range("M2").Value = 17.15870100078359
$ws.Range("N2").Value = 19.96945516764762
$ws.Range("B3").Value = 19.53062189691072
$ws.Range("C3").Value = 9.279205438660727
$ws.Range("D3").Value = 7.146797519895692
$ws.Range("E3").Value = 9.625757643671816
$ws.Range("F3").Value = 37.18173560654666
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 29.56658138265738
$ws.Range("L3").Value = 10.28598448907624
$ws.Range("M3").Value = 17.08375225915534
$ws.Range("N3").Value = 20.03970510188971
$ws.Range("B4").Value = 19.27467497848501
$ws.Range("C4").Value = 8.902594012467917
$ws.Range("D4").Value = 7.164480791490681
$ws.Range("E4").Value = 9.613020768527271
$ws.Range("F4").Value = 37.06904591238231
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 29.58179532446457
$ws.Range("L4").Value = 10.29436912991389
$ws.Range("M4").Value = 17.04140035155658
$ws.Range("N4").Value = 20.08481043669396
$ws.Range("B5").Value = 19.17083645641057
$ws.Range("C5").Value = 8.744893272948401
$ws.Range("D5").Value = 7.171885444953934
$ws.Range("E5").Value = 9.607763820629899
$ws.Range("F5").Value = 37.025898950411
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 29.58979024439698
$ws.Range("L5").Value = 10.29818027449676
$ws.Range("M5").Value = 17.02507607616613
$ws.Range("N5").Value = 20.10368820845638
$ws.Range("B6").Value = 19.15362674653124
$ws.Range("C6").Value = 8.718456267096631
$ws.Range("D6").Value = 7.1731269906108
$ws.Range("E6").Value = 9.606886885084551
$ws.Range("F6").Value = 37.01890283923579
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 29.59122603435729
$ws.Range("L6").Value = 10.29883692684895
$ws.Range("M6").Value = 17.02242222061761
$ws.Range("N6").Value = 20.10685290020146
$ws.Range("B7").Value = 19.27327249349795
$ws.Range("C7").Value = 8.900484128665443
$ws.Range("D7").Value = 7.164579848556964
$ws.Range("E7").Value = 9.612950141121276
$ws.Range("F7").Value = 37.06845274456721
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 29.58189588660913
$ws.Range("L7").Value = 10.29441893178899
$ws.Range("M7").Value = 17.04117639761869
$ws.Range("N7").Value = 20.08506301512488
$ws.Range("B8").Value = 19.80487072862391
$ws.Range("C8").Value = 9.666246758341975
$ws.Range("D8").Value = 7.128605157179305
$ws.Range("E8").Value = 9.639121801429003
$ws.Range("F8").Value = 37.3113863915969
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 29.55649881818186
$ws.Range("L8").Value = 10.27838531863757
$ws.Range("M8").Value = 17.1321054016085
$ws.Range("N8").Value = 19.99326887228443
$ws.Range("B9").Value = 20.85005655740153
$ws.Range("C9").Value = 11.04298256727099
$ws.Range("D9").Value = 7.064285177218053
$ws.Range("E9").Value = 9.688352488207929
$ws.Range("F9").Value = 37.8727963831178
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 29.56372884662326
$ws.Range("L9").Value = 10.25938689211187
$ws.Range("M9").Value = 17.33891387664097
$ws.Range("N9").Value = 19.82884798053212
$ws.Range("B10").Value = 21.61080840404817
$ws.Range("C10").Value = 12.02505508147887
$ws.Range("D10").Value = 7.020791527676787
$ws.Range("E10").Value = 9.723257075612759
$ws.Range("F10").Value = 38.33468505712664
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 29.60404877431933
$ws.Range("L10").Value = 10.25301170929248
$ws.Range("M10").Value = 17.5073530962333
$ws.Range("N10").Value = 19.71746922808123
$ws.Range("B11").Value = 21.95360272834502
$ws.Range("C11").Value = 12.44626122269932
$ws.Range("D11").Value = 7.001815042190892
$ws.Range("E11").Value = 9.738863542812561
$ws.Range("F11").Value = 38.55494336143999
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 29.63002401485448
$ws.Range("L11").Value = 10.25175595347936
$ws.Range("M11").Value = 17.58735445500433
$ws.Range("N11").Value = 19.66882873274604
$ws.Range("B12").Value = 22.08280422908948
$ws.Range("C12").Value = 12.60210883532569
$ws.Range("D12").Value = 6.994744973427484
$ws.Range("E12").Value = 9.744734325829187
$ws.Range("F12").Value = 38.63975215803095
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 29.64095854831385
$ws.Range("L12").Value = 10.25151654102807
$ws.Range("M12").Value = 17.61811555998526
$ws.Range("N12").Value = 19.65069999959346
$ws.Range("B13").Value = 22.05500746367288
$ws.Range("C13").Value = 12.56870647260056
$ws.Range("D13").Value = 6.996262490436296
$ws.Range("E13").Value = 9.743471683755631
$ws.Range("F13").Value = 38.62142565934115
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 29.63855475674452
$ws.Range("L13").Value = 10.25155760782636
$ws.Range("M13").Value = 17.61147018239839
$ws.Range("N13").Value = 19.65459144828034
$ws.Range("B14").Value = 21.96424513507431
$ws.Range("C14").Value = 12.45915602205463
$ws.Range("D14").Value = 7.001231062635514
$ws.Range("E14").Value = 9.739347318028562
$ws.Range("F14").Value = 38.56189287504408
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 29.63090159487043
$ws.Range("L14").Value = 10.25173152702725
$ws.Range("M14").Value = 17.58987595749124
$ws.Range("N14").Value = 19.66733145726879
$ws.Range("B15").Value = 21.90856760135151
$ws.Range("C15").Value = 12.39157791446598
$ws.Range("D15").Value = 7.00428953983904
$ws.Range("E15").Value = 9.736815940866451
$ws.Range("F15").Value = 38.52560813789233
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 29.62635683922865
$ws.Range("L15").Value = 10.25186879515734
$ws.Range("M15").Value = 17.57670901235885
$ws.Range("N15").Value = 19.67517286544607
$ws.Range("B16").Value = 21.58832824136612
$ws.Range("C16").Value = 11.99701590326107
$ws.Range("D16").Value = 7.022047928678587
$ws.Range("E16").Value = 9.722231690195063
$ws.Range("F16").Value = 38.32048997809498
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 29.60250491093582
$ws.Range("L16").Value = 10.25312683722682
$ws.Range("M16").Value = 17.50219114784126
$ws.Range("N16").Value = 19.72068868578458
$ws.Range("B17").Value = 21.39093328189137
$ws.Range("C17").Value = 11.74843618757464
$ws.Range("D17").Value = 7.033149026883363
$ws.Range("E17").Value = 9.713215137276832
$ws.Range("F17").Value = 38.19721508169564
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 29.58982839921149
$ws.Range("L17").Value = 10.25431954699086
$ws.Range("M17").Value = 17.45732840890877
$ws.Range("N17").Value = 19.74912941739431
$ws.Range("B18").Value = 21.27709625699113
$ws.Range("C18").Value = 11.60305428601886
$ws.Range("D18").Value = 7.039610255608826
$ws.Range("E18").Value = 9.708003468975665
$ws.Range("F18").Value = 38.12726892600399
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 29.58325573837214
$ws.Range("L18").Value = 10.25516036156645
$ws.Range("M18").Value = 17.43184366983732
$ws.Range("N18").Value = 19.7656785538697
$ws.Range("B19").Value = 21.23850564323893
$ws.Range("C19").Value = 11.55341675516557
$ws.Range("D19").Value = 7.041811011914248
$ws.Range("E19").Value = 9.706234492093198
$ws.Range("F19").Value = 38.10375262649178
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 29.58115372195601
$ws.Range("L19").Value = 10.25547164284378
$ws.Range("M19").Value = 17.42327036058393
$ws.Range("N19").Value = 19.77131460473426
$ws.Range("B20").Value = 21.41197844328621
$ws.Range("C20").Value = 11.77514673616277
$ws.Range("D20").Value = 7.031959415721318
$ws.Range("E20").Value = 9.714177612445662
$ws.Range("F20").Value = 38.21023912526167
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 29.59110346247987
$ws.Range("L20").Value = 10.25417656204046
$ws.Range("M20").Value = 17.46207123375724
$ws.Range("N20").Value = 19.74608211533769
$ws.Range("B21").Value = 21.9909217568994
$ws.Range("C21").Value = 12.49143266157373
$ws.Range("D21").Value = 6.999768529929234
$ws.Range("E21").Value = 9.740559803674145
$ws.Range("F21").Value = 38.57934152134627
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 29.6331197077568
$ws.Range("L21").Value = 10.2516740377813
$ws.Range("M21").Value = 17.59620620709446
$ws.Range("N21").Value = 19.66358153524597
$ws.Range("B22").Value = 22.36569777553004
$ws.Range("C22").Value = 12.93827738187672
$ws.Range("D22").Value = 6.979405377324428
$ws.Range("E22").Value = 9.757574507535207
$ws.Range("F22").Value = 38.82871335967624
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 29.66698085763245
$ws.Range("L22").Value = 10.25141453438906
$ws.Range("M22").Value = 17.68657974582427
$ws.Range("N22").Value = 19.61135466515732
$ws.Range("B23").Value = 22.16604472968869
$ws.Range("C23").Value = 12.70172963422835
$ws.Range("D23").Value = 6.990211906318808
$ws.Range("E23").Value = 9.748514252249324
$ws.Range("F23").Value = 38.69489344855867
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 29.64832290248382
$ws.Range("L23").Value = 10.25142726280296
$ws.Range("M23").Value = 17.63810445112109
$ws.Range("N23").Value = 19.63907463481921
$ws.Range("B24").Value = 21.40246501189595
$ws.Range("C24").Value = 11.76307859244121
$ws.Range("D24").Value = 7.032496992643223
$ws.Range("E24").Value = 9.713742564426331
$ws.Range("F24").Value = 38.20434806591619
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 29.59052477863699
$ws.Range("L24").Value = 10.25424072237943
$ws.Range("M24").Value = 17.45992604355642
$ws.Range("N24").Value = 19.74745918323168
$ws.Range("B25").Value = 20.56794662244408
$ws.Range("C25").Value = 10.67371753636325
$ws.Range("D25").Value = 7.081022366869359
$ws.Range("E25").Value = 9.675257170029759
$ws.Range("F25").Value = 37.71205026245645
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 29.55563854795122
$ws.Range("L25").Value = 10.26319430073802
$ws.Range("M25").Value = 17.28000198889912
$ws.Range("N25").Value = 19.87166793027212
